$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-37 down to 29-38
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with data (same structure as surrounding rows,
# preserving the date cell style used by column D)
$ws.Cells.Item(28, 1).Value = 4
$ws.Cells.Item(28, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(28, 3).Value = "Los Lagos"
$ws.Cells.Item(28, 4).Value = 44876
$ws.Cells.Item(28, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(28, 5).Value = 10
$ws.Cells.Item(28, 6).Value = 100112013
$ws.Cells.Item(28, 7).Value = "Alcachofa"
$ws.Cells.Item(28, 8).Value = "Española"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 200
$ws.Cells.Item(28, 11).Value = 12000
$ws.Cells.Item(28, 12).Value = 12000
$ws.Cells.Item(28, 13).Value = 12000
$ws.Cells.Item(28, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(28, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(28, 16).Value = 400
$ws.Cells.Item(28, 17).Value = 30
$ws.Cells.Item(28, 18).Value = "Hortaliza"
